# Update shopping cart & login: insert a "CreatDate" column into the
# Products sheet (between Available and CategoryId) and populate it
# with date values, formatted as yyyy-mm-dd, right-aligned.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# Insert a new column at F; this shifts the old CategoryId column (F) to G,
# carrying its data/styles along automatically.
$ws.Columns.Item(6).Insert()

# Header cell for the new column.
$ws.Range("F1").Value = "CreatDate"

# Date values (Excel serial date numbers) for rows 2..47, matching the
# CreatDate column of the Products sheet.
$dates = @(37592,38029,37021,36621,36960,37326,37633,38056,38634,38787,39212,39632,39875,40369,36979,40704,38903,41155,40612,37262,39423,40947,41681,44321,38238,40211,38787,39880,39817,44682,44502,44745,44442,44661,44449,44480,44537,44478,44389,44296,44267,44236,44199,44411,44540,44298)

$startRow = 2
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Range("F" + $row)
    $cell.HorizontalAlignment = -4152
    $cell.NumberFormat = "yyyy\-mm\-dd"
    $cell.Value = $dates[$i]
}
